$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 (time changed 20:03 -> 20:33)
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 20:33"

# Pakistan's case counts grew past Chile's, so the two countries swap rows
# (the table is sorted descending by "Casos totales", column B).
# Row 26 used to be Chile, row 27 used to be Pakistan; now row 26 is the
# (updated) Pakistan figures and row 27 holds Chile's (unchanged) figures.
$ws.Range("A26").Value = "Pakistan"
$ws.Range("A27").Value = "Chile"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1246462
$ws.Range("C4").Value = 8829
$ws.Range("D4").Value = 203753
$ws.Range("E4").Value = 969475
$ws.Range("G4").Value = 963
$ws.Range("H4").Value = 73234

# --- Row 8: Francia ---
$ws.Range("D8").Value = 53972
$ws.Range("E8").Value = 90770
$ws.Range("F8").Value = 3147
$ws.Range("G8").Value = 278
$ws.Range("H8").Value = 25809

# --- Row 12: Brasil ---
$ws.Range("B12").Value = 121600
$ws.Range("C12").Value = 6885
$ws.Range("E12").Value = 65357
$ws.Range("G12").Value = 101
$ws.Range("H12").Value = 8022

# --- Row 16: India ---
$ws.Range("B16").Value = 52559
$ws.Range("C16").Value = 3159
$ws.Range("D16").Value = 15257
$ws.Range("E16").Value = 35517
$ws.Range("G16").Value = 92
$ws.Range("H16").Value = 1785

# --- Row 26: Pakistan (new figures, now ranked above Chile) ---
$ws.Range("B26").Value = 23214
$ws.Range("C26").Value = 1165
$ws.Range("D26").Value = 6281
$ws.Range("E26").Value = 16389
$ws.Range("F26").Value = 111
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 544

# --- Row 27: Chile (figures unchanged, just shifted down one row) ---
$ws.Range("B27").Value = 23048
$ws.Range("C27").Value = 1032
$ws.Range("D27").Value = 11189
$ws.Range("E27").Value = 11578
$ws.Range("F27").Value = 470
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 281

# --- Row 73: Uzbekistan ---
$ws.Range("B73").Value = 2233
$ws.Range("C73").Value = 26
$ws.Range("E73").Value = 646

# --- Row 101: Sri Lanka ---
$ws.Range("B101").Value = 797
$ws.Range("C101").Value = 26
$ws.Range("E101").Value = 573
